$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (pure numeric-looking strings).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "28.042.33"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "1.914.99"
$ws.Range("E3").Value = "  -2.86%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.31%  "
$ws.Range("D5").Value = "329.20"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("E7").Value = "  -5.61%  "
$ws.Range("E8").Value = "  -4.27%  "
$ws.Range("D9").Value = "53.03"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").Value = "0.08380"
$ws.Range("E10").Value = "  -10.06%  "
$ws.Range("E11").Value = "  -5.06%  "
$ws.Range("D12").Value = "22.07"
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").Value = "7.445"
$ws.Range("E13").Value = "  -5.48%  "
$ws.Range("D14").Value = "1.872.00"
$ws.Range("E14").Value = "  -5.17%  "
$ws.Range("D15").Value = "6.054"
$ws.Range("E15").Value = "  -6.16%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "89.44"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("E18").Value = "  -4.30%  "
$ws.Range("D19").Value = "0.06584"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").Value = "17.95"
$ws.Range("E20").Value = "  -6.21%  "
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "5.704"
$ws.Range("E22").Value = "  -4.06%  "
$ws.Range("D23").Value = "28.009.56"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("E24").Value = "  -5.17%  "
$ws.Range("D25").Value = "2.286"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "2.098.15"
$ws.Range("E26").Value = "  -5.39%  "
$ws.Range("D27").Value = "153.93"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").Value = "19.94"
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("D29").Value = "2.127"
$ws.Range("E29").Value = "  -6.07%  "
$ws.Range("D30").Value = "5.685"
$ws.Range("E30").Value = "  -8.71%  "
$ws.Range("E31").Value = "  -3.51%  "
$ws.Range("D32").Value = "0.9698"
$ws.Range("E32").Value = "  -7.14%  "
$ws.Range("D33").Value = "0.09550"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").Value = "1.440"
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").Value = "3.641"
$ws.Range("E35").Value = "  -2.90%  "
$ws.Range("D36").Value = "5.519"
$ws.Range("E36").Value = "  -4.88%  "
$ws.Range("D37").Value = "8.818"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("D38").Value = "0.02300"
$ws.Range("E38").Value = "  -4.69%  "
$ws.Range("D39").Value = "0.06137"
$ws.Range("E39").Value = "  -4.07%  "
$ws.Range("E40").Value = "  -8.46%  "
$ws.Range("D41").Value = "0.6117"
$ws.Range("E41").Value = "  -5.41%  "
$ws.Range("D42").Value = "11.00"
$ws.Range("E42").Value = "  -4.29%  "
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("E44").Value = "  -5.18%  "
$ws.Range("D45").Value = "1.301"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("D46").Value = "0.5848"
$ws.Range("E46").Value = "  -5.48%  "
$ws.Range("D47").Value = "12.74"
$ws.Range("E47").Value = "  -4.40%  "
$ws.Range("E48").Value = "  -7.37%  "
$ws.Range("D49").Value = "3.463"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").Value = "0.06824"
$ws.Range("E50").Value = "  -1.98%  "
